# Update cryptos list prices/volume percentages (and restore WEMIXToken/BabyDogeCoin row order)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.123.38"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.791.42"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'228.62"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'32.72"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").Value = "'0.0715"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "2.049.02"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "'11.13"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "1.800.00"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "'0.626"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "34.063.05"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "'4.17"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").Value = "'68.52"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "'245.62"
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("D20").Value = "0.0₃0790"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'10.78"
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("D23").Value = "'4.11"
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("E24").Value = "  -2.98%  "
$ws.Range("D25").Value = "'160.63"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").Value = "'16.37"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("D35").Value = "1.401.17"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("D36").Value = "'0.661"
$ws.Range("E36").Value = "  +3.32%  "
$ws.Range("D37").Value = "'1.05"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").Value = "'2.23"
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").Value = "'0.920"
$ws.Range("E41").Value = "  -4.30%  "
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("D44").Value = "'13.32"
$ws.Range("E44").Value = "  +11.79%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.08"
$ws.Range("E45").Value = "  +3.35%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0138"
$ws.Range("E46").Value = "  +9.94%  "
$ws.Range("D47").Value = "'109.40"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("D48").Value = "'0.0498"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "'5.84"
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("D50").Value = "1.948.34"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("E51").Value = "  +0.24%  "
